$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_LT"
$ws.Range("A2").Value = 57.869742198100404
$ws.Range("A3").Value = 58.208955223880601
$ws.Range("A4").Value = 58.412483039348707
$ws.Range("A5").Value = 57.394843962008146
$ws.Range("A6").Value = 60.651289009497965
$ws.Range("A7").Value = 60.651289009497965
$ws.Range("A8").Value = 55.223880597014926
$ws.Range("A9").Value = 56.580732700135684
$ws.Range("A10").Value = 57.055630936227949
$ws.Range("A11").Value = 55.970149253731336
$ws.Range("A12").Value = 56.919945725915881
$ws.Range("A13").Value = 56.309362279511532
$ws.Range("A14").Value = 53.79918588873813
$ws.Range("A15").Value = 53.188602442333789
$ws.Range("A16").Value = 54.613297150610585
$ws.Range("A17").Value = 55.156037991858888
$ws.Range("A18").Value = 57.191316146540025
$ws.Range("A19").Value = 57.055630936227949
$ws.Range("A20").Value = 58.887381275440973
$ws.Range("A21").Value = 59.905020352781548
$ws.Range("A22").Value = 60.040705563093624
$ws.Range("A23").Value = 59.362279511533245
$ws.Range("A24").Value = 60.447761194029844
$ws.Range("A25").Value = 60.719131614653996
$ws.Range("A26").Value = 59.633649932157397
$ws.Range("A27").Value = 57.191316146540025
$ws.Range("A28").Value = 58.00542740841248
$ws.Range("A29").Value = 53.45997286295794
$ws.Range("A30").Value = 53.324287652645864
$ws.Range("A31").Value = 51.017639077340569
$ws.Range("A32").Value = 62.550881953867034
$ws.Range("A33").Value = 60.854816824966072
$ws.Range("A34").Value = 60.854816824966072
$ws.Range("A35").Value = 59.497964721845321
$ws.Range("A36").Value = 60.58344640434192
$ws.Range("A37").Value = 52.51017639077341
$ws.Range("A38").Value = 61.668928086838534
$ws.Range("A39").Value = 58.412483039348707
$ws.Range("A40").Value = 58.208955223880601
$ws.Range("A41").Value = 57.327001356852101
$ws.Range("A42").Value = 53.663500678426054
$ws.Range("A43").Value = 54.274084124830388
$ws.Range("A44").Value = 54.816824966078691
$ws.Range("A45").Value = 56.784260515603805
$ws.Range("A46").Value = 56.648575305291729
$ws.Range("A47").Value = 55.495251017639077
$ws.Range("A48").Value = 53.120759837177744
$ws.Range("A49").Value = 56.71641791044776
